# "updated boms with antenna"
# Update the mccoy_bom sheet:
#  - swap the L1/L2 inductor part for a higher-current one (row 8)
#  - add three new BOM lines: U3 (1.1V LDO), U2 (1.8V LDO), and the SMA antenna connector
#  - move the active selection to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mccoy_bom")

# --- Row 8: L1, L2 inductor changed from MLF1005G2R2JT000 to MLZ1005M2R2WT000 ---
$ws.Range("B8").Value = "FIXED IND 2.2UH 350MA 550MOHM SM"
$ws.Range("H8").Value = "https://www.digikey.com/en/products/detail/tdk-corporation/MLZ1005M2R2WT000/2465140"

# --- Row 10: new line item - U3, 1.1V LDO regulator ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "IC REG LINEAR 1.1V 300MA SOT23-5"
$ws.Range("C10").Value = "U3"
$ws.Range("F10").Value = "Digi-Key"
$ws.Range("H10").Value = "https://www.digikey.com/en/products/detail/texas-instruments/TLV70311DBVR/7776390"
$ws.Range("I10").Value = 1

# --- Row 11: new line item - U2, 1.8V LDO regulator ---
$ws.Range("B11").Value = "IC REG LINEAR 1.8V 250MA SOT23-5"
$ws.Range("C11").Value = "U2"
$ws.Range("H11").Value = "https://www.digikey.com/en/products/detail/texas-instruments/LP5907MFX-1-8-NOPB/3911201"
$ws.Range("I11").Value = 1

# --- Row 12: new line item - SMA antenna connector ---
$ws.Range("B12").Value = "CONN SMA RCPT STR 50 OHM PCB"
$ws.Range("H12").Value = "https://www.digikey.com/en/products/detail/molex/0733910060/1465165"
$ws.Range("I12").Value = 1

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()
